$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell O1 "specific_prompt", copying formatting from N1 (the previous last header)
$ws.Range("N1").Copy() | Out-Null
$ws.Range("O1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("O1").Value = "specific_prompt"

# Update A2 value: 25 -> 30
$ws.Range("A2").Value = 30

# Update E2 value: email changed
$ws.Range("E2").Value = "charlie.brown@example.com"

# Update M2 note text (new timestamp + content)
$ws.Range("M2").Value = "`n[2025-09-24 17:43:24] The customer is from the UK, likes football, and has dust allergies, but the conversation didn't go further as the customer didn't express any specific requirements or interests beyond that."

# Keep row 2 height auto (setting long text above can trigger an explicit row height); re-autofit to avoid a stored custom height
$ws.Rows.Item(2).AutoFit() | Out-Null

# Clear N2 (tasks) value - becomes blank
$ws.Range("N2").Value = ""

# Set new O2 value (specific_prompt content for row 2)
$ws.Range("O2").Value = "he is getting alergies of dust particals"
